$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit re-shuffles the per-occurrence data (Id, Antal, Ost, Nord — columns
# A, I, Q, R) among rows 35-47 of the sightings table; every other column in
# those rows is unchanged. Column I ("Antal") is stored as text, so its
# value is staged through a scratch column (via Range.Copy, which preserves
# the text type) rather than via .Value (which would coerce a numeric-looking
# string back into a Number).

$scratchCol = "ZZ"

# Row N's new content comes from row Source(N)'s old content.
$rowSource = @{
    35 = 46
    36 = 45
    37 = 47
    38 = 37
    40 = 36
    41 = 38
    42 = 40
    43 = 41
    44 = 35
    45 = 43
    46 = 44
    47 = 42
}

# 1) Stage the current (pre-edit) "Antal" (I) text for every row that is used
#    as a source, into scratch cells keyed by destination row, before any
#    destination cell gets overwritten.
foreach ($destRow in $rowSource.Keys) {
    $srcRow = $rowSource[$destRow]
    $ws.Range("I$srcRow").Copy($ws.Range("$scratchCol$destRow"))
}

# 2) New Id (A) / Ost (Q) / Nord (R) numeric values, per destination row.
$newA = @{
    35 = 111675584
    36 = 111675582
    37 = 111675575
    38 = 111675574
    40 = 111675580
    41 = 111675579
    42 = 111675572
    43 = 111675577
    44 = 111675573
    45 = 111675571
    46 = 111675583
    47 = 111675578
}
$newQ = @{
    35 = 690414.984509701
    36 = 690352.3333891984
    37 = 690480.7418955797
    38 = 690486.6986671695
    40 = 690370.5537696742
    41 = 690425.8424831247
    42 = 690494.5947179901
    43 = 690430.9193086301
    44 = 690487.9917822112
    45 = 690509.4285896254
    46 = 690415.8809986882
    47 = 690368.3990222017
}
$newR = @{
    35 = 6661422.355185229
    36 = 6661470.655078794
    37 = 6661091.463633558
    38 = 6661102.281881573
    40 = 6661292.946251329
    41 = 6661357.862056008
    42 = 6661104.692649405
    43 = 6661356.623615522
    44 = 6661106.352564453
    45 = 6661040.900344189
    46 = 6661424.403280765
    47 = 6661295.837351476
}

foreach ($row in $newA.Keys) {
    $ws.Range("A$row").Value = $newA[$row]
    $ws.Range("Q$row").Value = $newQ[$row]
    $ws.Range("R$row").Value = $newR[$row]
}

# 3) Rows whose "Antal" (I) text actually changes value.
$rowsWithNewI = @(36, 38, 40, 41, 42, 43, 44, 45, 46)
foreach ($row in $rowsWithNewI) {
    $ws.Range("$scratchCol$row").Copy($ws.Range("I$row"))
}

# 4) Clean up the scratch column.
foreach ($destRow in $rowSource.Keys) {
    $ws.Range("$scratchCol$destRow").ClearContents()
}
